# ADD results from server
# Updates the row-2 result values on each year sheet (2025..2050) to the
# freshly computed figures returned by the server re-run.

$wb = $excel.ActiveWorkbook

$updates = @{
    "2025" = @{
        A = 4386.58098096851
        B = 1887.750563722889
        E = 13582.352998632
        G = 4231.516049511674
        H = 49322.36395174918
        I = 40307.21206296
        N = 5169.61678787112
        O = 6937.012684477218
    }
    "2030" = @{
        A = 9530.678248741482
        B = 15047.9874167341
        E = 24805.53651836642
        G = 4231.516049511674
        H = 59837.73428632267
        I = 73842.81869365374
        N = 7603.607993347517
        O = 12201.72900914311
    }
    "2035" = @{
        A = 18114.08743212745
        B = 15048.22929088937
        E = 24805.53651836642
        G = 4231.516049511674
        H = 59837.73428632267
        I = 91505.69078165847
        M = 2812.883060279026
        N = 8803.659275398724
        O = 17930.2177511366
    }
    "2040" = @{
        A = 18114.08743212745
        B = 15048.22929088937
        E = 24805.53651836642
        G = 4231.516049511674
        H = 59837.73428632267
        I = 91505.69078165847
        M = 2812.883060279026
        N = 8803.659275398724
        O = 17930.2177511366
    }
    "2045" = @{
        A = 18114.08743212745
        B = 15048.22929088937
        E = 24805.53651836642
        G = 4231.516049511674
        H = 59837.73428632267
        I = 91505.69078165847
        M = 2812.883060279026
        N = 8803.659275398724
        O = 17930.2177511366
    }
    "2050" = @{
        A = 18114.08743212745
        B = 15048.22929088937
        E = 24805.53651836642
        G = 4231.516049511674
        H = 59837.73428632267
        I = 91505.69078165847
        M = 2812.883060279026
        N = 8803.659275398724
        O = 17930.2177511366
    }
}

foreach ($sheetName in $updates.Keys) {
    $sheetNameStr = [string]$sheetName
    $ws = $wb.Worksheets.Item($sheetNameStr)
    $cellValues = $updates[$sheetName]
    foreach ($col in $cellValues.Keys) {
        $colStr = [string]$col
        $ws.Range("$colStr`2").Value = $cellValues[$col]
    }
}
